$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.21398521684939
$ws.Range("C2").Value = 12.31560227565291
$ws.Range("E2").Value = 13.36223393530197
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 24.27843624243556
$ws.Range("H2").Value = 13.18575121505748
$ws.Range("L2").Value = 9.814528489030314
$ws.Range("M2").Value = 13.83754049242213
$ws.Range("N2").Value = 17.57552532045061
$ws.Range("O2").Value = 19.46715587116932

# Row 3
$ws.Range("B3").Value = 12.76271241188098
$ws.Range("C3").Value = 12.24719418083618
$ws.Range("E3").Value = 13.41174372781251
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 24.3244957366343
$ws.Range("H3").Value = 13.23195892781406
$ws.Range("L3").Value = 9.820490745795702
$ws.Range("M3").Value = 13.74776926705693
$ws.Range("N3").Value = 17.61480324029826
$ws.Range("O3").Value = 19.53824825278387

# Row 4
$ws.Range("B4").Value = 12.47861274572599
$ws.Range("C4").Value = 12.20523890311344
$ws.Range("E4").Value = 13.44394344324337
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 24.36206271122449
$ws.Range("H4").Value = 13.26257631663888
$ws.Range("L4").Value = 9.825457767990565
$ws.Range("M4").Value = 13.6941395386984
$ws.Range("N4").Value = 17.6407056856032
$ws.Range("O4").Value = 19.58653506266554

# Row 5
$ws.Range("B5").Value = 12.3612485512356
$ws.Range("C5").Value = 12.18816241285769
$ws.Range("E5").Value = 13.45751866368636
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 24.37969638943608
$ws.Range("H5").Value = 13.27561767765303
$ws.Range("L5").Value = 9.82781096485796
$ws.Range("M5").Value = 13.67267685792054
$ws.Range("N5").Value = 17.65171089283085
$ws.Range("O5").Value = 19.60737515503186

# Row 6
$ws.Range("B6").Value = 12.34166957351984
$ws.Range("C6").Value = 12.18532837167309
$ws.Range("E6").Value = 13.45980023788753
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 24.38276453313254
$ws.Range("H6").Value = 13.27781727698055
$ws.Range("L6").Value = 9.828221605492161
$ws.Range("M6").Value = 13.66913715406968
$ws.Range("N6").Value = 17.65356548471367
$ws.Range("O6").Value = 19.61090578215668

# Row 7
$ws.Range("B7").Value = 12.47703613392674
$ws.Range("C7").Value = 12.20500850870697
$ws.Range("E7").Value = 13.44412468579953
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 24.36229112670279
$ws.Range("H7").Value = 13.26274991146565
$ws.Range("L7").Value = 9.825488170726777
$ws.Range("M7").Value = 13.6938484765852
$ws.Range("N7").Value = 17.64085228367224
$ws.Range("O7").Value = 19.58681141508862

# Row 8
$ws.Range("B8").Value = 13.05994153752638
$ws.Range("C8").Value = 12.2920105115075
$ws.Range("E8").Value = 13.37893171336416
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 24.29238446508953
$ws.Range("H8").Value = 13.20121744874455
$ws.Range("L8").Value = 9.816313586184357
$ws.Range("M8").Value = 13.80628854702068
$ws.Range("N8").Value = 17.58869817737344
$ws.Range("O8").Value = 19.4907046584724

# Row 9
$ws.Range("B9").Value = 14.14053950332914
$ws.Range("C9").Value = 12.46261686365514
$ws.Range("E9").Value = 13.26533709516214
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 24.22936851480515
$ws.Range("H9").Value = 13.0983780670335
$ws.Range("L9").Value = 9.808655789723632
$ws.Range("M9").Value = 14.03782983053735
$ws.Range("N9").Value = 17.50055991728384
$ws.Range("O9").Value = 19.33914790159121

# Row 10
$ws.Range("B10").Value = 14.88833202440489
$ws.Range("C10").Value = 12.5874465870851
$ws.Range("E10").Value = 13.19051291665977
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 24.22865031274865
$ws.Range("H10").Value = 13.03369863172507
$ws.Range("L10").Value = 9.809285565251896
$ws.Range("M10").Value = 14.21360734638015
$ws.Range("N10").Value = 17.44437833119457
$ws.Range("O10").Value = 19.25046517787338

# Row 11
$ws.Range("B11").Value = 15.21710041074613
$ws.Range("C11").Value = 12.64401337596171
$ws.Range("E11").Value = 13.15833718982454
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 24.23827541051871
$ws.Range("H11").Value = 13.00663821267876
$ws.Range("L11").Value = 9.810919712082523
$ws.Range("M11").Value = 14.29456792428026
$ws.Range("N11").Value = 17.42067242497625
$ws.Range("O11").Value = 19.21507606672535

# Row 12
$ws.Range("B12").Value = 15.3398569397975
$ws.Range("C12").Value = 12.66539291744417
$ws.Range("E12").Value = 13.14642005910953
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 24.24335277880733
$ws.Range("H12").Value = 12.9967310986754
$ws.Range("L12").Value = 9.811731250167089
$ws.Range("M12").Value = 14.32534899344218
$ws.Range("N12").Value = 17.41196115280299
$ws.Range("O12").Value = 19.20238993191358

# Row 13
$ws.Range("B13").Value = 15.31349802745574
$ws.Range("C13").Value = 12.66079043004004
$ws.Range("E13").Value = 13.14897475659869
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 24.24219555987991
$ws.Range("H13").Value = 12.99884964236169
$ws.Range("L13").Value = 9.811547917448499
$ws.Range("M13").Value = 14.31871459464658
$ws.Range("N13").Value = 17.41382547593837
$ws.Range("O13").Value = 19.20509027878429

# Row 14
$ws.Range("B14").Value = 15.22723508858589
$ws.Range("C14").Value = 12.64577316179452
$ws.Range("E14").Value = 13.15735141057348
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 24.23866441633817
$ws.Range("H14").Value = 13.00581633078745
$ws.Range("L14").Value = 9.810982621785609
$ws.Range("M14").Value = 14.29709795316363
$ws.Range("N14").Value = 17.41995042286887
$ws.Range("O14").Value = 19.21401802707773

# Row 15
$ws.Range("B15").Value = 15.17416704744559
$ws.Range("C15").Value = 12.63656901780802
$ws.Range("E15").Value = 13.16251712044012
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 24.23668806168712
$ws.Range("H15").Value = 13.01012792715202
$ws.Range("L15").Value = 9.810661426851359
$ws.Range("M15").Value = 14.28387255016388
$ws.Range("N15").Value = 17.42373670667936
$ws.Range("O15").Value = 19.21957971650207

# Row 16
$ws.Range("B16").Value = 14.8666078228649
$ws.Range("C16").Value = 12.5837445890422
$ws.Range("E16").Value = 13.19265309102293
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 24.22822171667922
$ws.Range("H16").Value = 13.03551466426839
$ws.Range("L16").Value = 9.809205798684188
$ws.Range("M16").Value = 14.2083347229588
$ws.Range("N16").Value = 17.44596477205883
$ws.Range("O16").Value = 19.25287784298846

# Row 17
$ws.Range("B17").Value = 14.67493317785878
$ws.Range("C17").Value = 12.55127561852067
$ws.Range("E17").Value = 13.21161700549914
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 24.22557844696193
$ws.Range("H17").Value = 13.05169396321079
$ws.Range("L17").Value = 9.808657250535841
$ws.Range("M17").Value = 14.16223589391508
$ws.Range("N17").Value = 17.46007473244137
$ws.Range("O17").Value = 19.27457574250419

# Row 18
$ws.Range("B18").Value = 14.5636185879808
$ws.Range("C18").Value = 12.53258029747523
$ws.Range("E18").Value = 13.22269986203801
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 24.22499484915824
$ws.Range("H18").Value = 13.06122222478703
$ws.Range("L18").Value = 9.808468657311735
$ws.Range("M18").Value = 14.13581629029293
$ws.Range("N18").Value = 17.46836470336434
$ws.Range("O18").Value = 19.28752185721084

# Row 19
$ws.Range("B19").Value = 14.52574915371981
$ws.Range("C19").Value = 12.5262472359152
$ws.Range("E19").Value = 13.22648245756425
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 24.22495806275501
$ws.Range("H19").Value = 13.06448650996171
$ws.Range("L19").Value = 9.808426632286569
$ws.Range("M19").Value = 14.12688804438258
$ws.Range("N19").Value = 17.47120149843171
$ws.Range("O19").Value = 19.29198515021934

# Row 20
$ws.Range("B20").Value = 14.69544861549012
$ws.Range("C20").Value = 12.55473413827374
$ws.Range("E20").Value = 13.20958012485423
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 24.22576286153148
$ws.Range("H20").Value = 13.04994863363529
$ws.Range("L20").Value = 9.808702516464066
$ws.Range("M20").Value = 14.16713348595636
$ws.Range("N20").Value = 17.45855466847874
$ws.Range("O20").Value = 19.27221770888788

# Row 21
$ws.Range("B21").Value = 15.25262059523088
$ws.Range("C21").Value = 12.65018528862993
$ws.Range("E21").Value = 13.15488373990267
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 24.23966271719416
$ws.Range("H21").Value = 13.00376081313454
$ws.Range("L21").Value = 9.81114344132847
$ws.Range("M21").Value = 14.30344410927067
$ws.Range("N21").Value = 17.41814417270843
$ws.Range("O21").Value = 19.21137630653515

# Row 22
$ws.Range("B22").Value = 15.60657893263462
$ws.Range("C22").Value = 12.71232429021192
$ws.Range("E22").Value = 13.1206931671416
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 24.25709644223095
$ws.Range("H22").Value = 12.97555680375682
$ws.Range("L22").Value = 9.813861505356764
$ws.Range("M22").Value = 14.39323865217176
$ws.Range("N22").Value = 17.39328166654863
$ws.Range("O22").Value = 19.17578137928851

# Row 23
$ws.Range("B23").Value = 15.41862687546564
$ws.Range("C23").Value = 12.67918501208086
$ws.Range("E23").Value = 13.13879909725006
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 24.24702776376817
$ws.Range("H23").Value = 12.99042831142848
$ws.Range("L23").Value = 9.812308463489606
$ws.Range("M23").Value = 14.34525567190939
$ws.Range("N23").Value = 17.40640979465766
$ws.Range("O23").Value = 19.19439679355301

# Row 24
$ws.Range("B24").Value = 14.68617706155363
$ws.Range("C24").Value = 12.55317062947214
$ws.Range("E24").Value = 13.2105004370898
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 24.2256765721192
$ws.Range("H24").Value = 13.05073699141731
$ws.Range("L24").Value = 9.808681656772187
$ws.Range("M24").Value = 14.16491902378326
$ws.Range("N24").Value = 17.4592413349959
$ws.Range("O24").Value = 19.27328230676798

# Row 25
$ws.Range("B25").Value = 13.85578519712776
$ws.Range("C25").Value = 12.41651643743894
$ws.Range("E25").Value = 13.29454752056018
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 24.23843571245483
$ws.Range("H25").Value = 13.12428949792556
$ws.Range("L25").Value = 9.809625349905863
$ws.Range("M25").Value = 13.97411807896242
$ws.Range("N25").Value = 17.52289490620479
$ws.Range("O25").Value = 19.37617916808807
